# Auto-generated script to apply numeric updates described in the commit diff.
# Updates currentAveragePrice / Leve price / profit columns (H, I, J, K, L, M, N)
# across several worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 10418124
$ws.Range("I15").Value = 10418124
$ws.Range("K15").Value = 31254372
$ws.Range("M15").Value = -31254203
$ws.Range("H20").Value = 3950
$ws.Range("I20").Value = 3950
$ws.Range("K20").Value = 3950
$ws.Range("M20").Value = -3720
$ws.Range("H35").Value = 3950
$ws.Range("I35").Value = 3950
$ws.Range("K35").Value = 3950
$ws.Range("M35").Value = -3571
$ws.Range("H76").Value = 3013.1428
$ws.Range("I76").Value = 2098.5
$ws.Range("J76").Value = 3379
$ws.Range("K76").Value = 2098.5
$ws.Range("L76").Value = 3379
$ws.Range("M76").Value = -1783.5
$ws.Range("N76").Value = -4009
$ws.Range("H79").Value = 3013.1428
$ws.Range("I79").Value = 2098.5
$ws.Range("J79").Value = 3379
$ws.Range("K79").Value = 2098.5
$ws.Range("L79").Value = 3379
$ws.Range("M79").Value = -1006.5
$ws.Range("N79").Value = -5563
$ws.Range("H80").Value = 990.53845
$ws.Range("I80").Value = 569.5714
$ws.Range("J80").Value = 1481.6666
$ws.Range("K80").Value = 1708.7142
$ws.Range("L80").Value = 4444.9998
$ws.Range("M80").Value = -710.7142000000001
$ws.Range("N80").Value = -6440.9998
$ws.Range("H83").Value = 990.53845
$ws.Range("I83").Value = 569.5714
$ws.Range("J83").Value = 1481.6666
$ws.Range("K83").Value = 5126.1426
$ws.Range("L83").Value = 13334.9994
$ws.Range("M83").Value = -134.1426000000001
$ws.Range("N83").Value = -23318.9994
$ws.Range("H138").Value = 4632.7466
$ws.Range("I138").Value = 1515.5
$ws.Range("J138").Value = 5004.955
$ws.Range("K138").Value = 4546.5
$ws.Range("L138").Value = 15014.865
$ws.Range("M138").Value = 593.5
$ws.Range("N138").Value = -25294.865
$ws.Range("H140").Value = 62882.25
$ws.Range("J140").Value = 61764.145
$ws.Range("L140").Value = 61764.145
$ws.Range("N140").Value = -72124.14499999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1344025.8
$ws.Range("I2").Value = 2909641.8
$ws.Range("J2").Value = 2069
$ws.Range("K2").Value = 2909641.8
$ws.Range("L2").Value = 2069
$ws.Range("M2").Value = -2909528.8
$ws.Range("N2").Value = -2295
$ws.Range("H45").Value = 58300.85
$ws.Range("I45").Value = 67707.06
$ws.Range("J45").Value = 4999
$ws.Range("K45").Value = 67707.06
$ws.Range("L45").Value = 4999
$ws.Range("M45").Value = -67330.06
$ws.Range("N45").Value = -5753
$ws.Range("H63").Value = 1959.2
$ws.Range("I63").Value = 2198.6667
$ws.Range("J63").Value = 1600
$ws.Range("K63").Value = 2198.6667
$ws.Range("L63").Value = 1600
$ws.Range("M63").Value = -1512.6667
$ws.Range("N63").Value = -2972
$ws.Range("H66").Value = 1959.2
$ws.Range("I66").Value = 2198.6667
$ws.Range("J66").Value = 1600
$ws.Range("K66").Value = 10993.3335
$ws.Range("L66").Value = 8000
$ws.Range("M66").Value = -7561.333500000001
$ws.Range("N66").Value = -14864
$ws.Range("H116").Value = 1344025.8
$ws.Range("I116").Value = 2909641.8
$ws.Range("J116").Value = 2069
$ws.Range("K116").Value = 2909641.8
$ws.Range("L116").Value = 2069
$ws.Range("M116").Value = -2907347.8
$ws.Range("N116").Value = -6657

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1344025.8
$ws.Range("I3").Value = 2909641.8
$ws.Range("J3").Value = 2069
$ws.Range("K3").Value = 2909641.8
$ws.Range("L3").Value = 2069
$ws.Range("M3").Value = -2909527.8
$ws.Range("N3").Value = -2297

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 198166.5
$ws.Range("J20").Value = 198166.5
$ws.Range("L20").Value = 198166.5
$ws.Range("N20").Value = -198638.5
$ws.Range("H30").Value = 198166.5
$ws.Range("J30").Value = 198166.5
$ws.Range("L30").Value = 198166.5
$ws.Range("N30").Value = -198348.5
$ws.Range("H31").Value = 5553.915
$ws.Range("I31").Value = 1981.9
$ws.Range("J31").Value = 6519.324
$ws.Range("K31").Value = 1981.9
$ws.Range("L31").Value = 6519.324
$ws.Range("M31").Value = -1686.9
$ws.Range("N31").Value = -7109.324
$ws.Range("H34").Value = 5553.915
$ws.Range("I34").Value = 1981.9
$ws.Range("J34").Value = 6519.324
$ws.Range("K34").Value = 1981.9
$ws.Range("L34").Value = 6519.324
$ws.Range("M34").Value = -1779.9
$ws.Range("N34").Value = -6923.324
$ws.Range("H62").Value = 41038.125
$ws.Range("J62").Value = 53733.332
$ws.Range("L62").Value = 53733.332
$ws.Range("N62").Value = -54981.332
$ws.Range("H65").Value = 41038.125
$ws.Range("J65").Value = 53733.332
$ws.Range("L65").Value = 268666.66
$ws.Range("N65").Value = -274906.66
$ws.Range("H122").Value = 2915.276
$ws.Range("I122").Value = 1805.5264
$ws.Range("J122").Value = 5023.8
$ws.Range("K122").Value = 5416.5792
$ws.Range("L122").Value = 15071.4
$ws.Range("M122").Value = -2966.5792
$ws.Range("N122").Value = -19971.4
$ws.Range("H128").Value = 198166.5
$ws.Range("J128").Value = 198166.5
$ws.Range("L128").Value = 198166.5
$ws.Range("N128").Value = -208126.5
$ws.Range("H134").Value = 1991.7
$ws.Range("I134").Value = 1990.3462
$ws.Range("K134").Value = 5971.0386
$ws.Range("M134").Value = -3436.0386
$ws.Range("H141").Value = 103972.914
$ws.Range("J141").Value = 109488.73
$ws.Range("L141").Value = 109488.73
$ws.Range("N141").Value = -119848.73

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 71.25
$ws.Range("I8").Value = 71.25
$ws.Range("K8").Value = 213.75
$ws.Range("M8").Value = -74.75
$ws.Range("H69").Value = 4700
$ws.Range("J69").Value = 4700
$ws.Range("L69").Value = 14100
$ws.Range("N69").Value = -15722
$ws.Range("H72").Value = 4700
$ws.Range("J72").Value = 4700
$ws.Range("L72").Value = 42300
$ws.Range("N72").Value = -50412
$ws.Range("H86").Value = 116.545456
$ws.Range("I86").Value = 282
$ws.Range("K86").Value = 846
$ws.Range("M86").Value = 340
$ws.Range("H89").Value = 116.545456
$ws.Range("I89").Value = 282
$ws.Range("K89").Value = 2538
$ws.Range("M89").Value = 3390
$ws.Range("H98").Value = 905.06665
$ws.Range("J98").Value = 953.6429000000001
$ws.Range("L98").Value = 2860.9287
$ws.Range("N98").Value = -5856.9287

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1138
$ws.Range("I113").Value = 1126.6
$ws.Range("J113").Value = 1195
$ws.Range("K113").Value = 1126.6
$ws.Range("L113").Value = 1195
$ws.Range("M113").Value = 1043.4
$ws.Range("N113").Value = -5535
$ws.Range("H140").Value = 79899
$ws.Range("J140").Value = 79899
$ws.Range("L140").Value = 79899
$ws.Range("N140").Value = -90259

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5824.25
$ws.Range("I132").Value = 5382.8184
$ws.Range("K132").Value = 16148.4552
$ws.Range("M132").Value = -13618.4552
$ws.Range("H139").Value = 78715
$ws.Range("J139").Value = 78715
$ws.Range("L139").Value = 78715
$ws.Range("N139").Value = -88995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 26250
$ws.Range("J48").Value = 26250
$ws.Range("L48").Value = 26250
$ws.Range("N48").Value = -27388
$ws.Range("H81").Value = 2988524.2
$ws.Range("I81").Value = 5209583
$ws.Range("K81").Value = 10419166
$ws.Range("M81").Value = -10418105
$ws.Range("H84").Value = 2988524.2
$ws.Range("I84").Value = 5209583
$ws.Range("K84").Value = 52095830
$ws.Range("M84").Value = -52090526
$ws.Range("H107").Value = 2935.48
$ws.Range("I107").Value = 3343.5334
$ws.Range("J107").Value = 2323.4
$ws.Range("K107").Value = 10030.6002
$ws.Range("L107").Value = 6970.200000000001
$ws.Range("M107").Value = -8110.600199999999
$ws.Range("N107").Value = -10810.2
$ws.Range("H113").Value = 942.8823
$ws.Range("I113").Value = 907
$ws.Range("J113").Value = 1008.6667
$ws.Range("K113").Value = 2721
$ws.Range("L113").Value = 3026.0001
$ws.Range("M113").Value = -551
$ws.Range("N113").Value = -7366.0001
$ws.Range("H133").Value = 62595.25
$ws.Range("J133").Value = 62595.25
$ws.Range("L133").Value = 62595.25
$ws.Range("N133").Value = -72715.25
